$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking": Right count 5 -> 4, Wrong penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total 85 -> 68, and summary text "85 / 140" -> "68 / 112"
$ws.Range("B12").Value = 68
$ws.Range("E12").Value = "68 / 112"
